# Generate Report for Handoff
# Updates the localization-status report:
#  - Sets "Priority" (column E) to "ht" for the newly-handed-off rows
#    on both the zh-cn and de-de sheets.
#  - Refreshes the "Latest Handoff Datetime" (column H on zh-cn/de-de)
#    and "Latest HO Xliff Generate Date" (column G on Overview) timestamps
#    for those same rows.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

# Overview sheet: refresh "Latest HO Xliff Generate Date" (column G)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-17 18:21:01"
}

# zh-cn sheet: set Priority (column E) and Latest Handoff Datetime (column H)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-17 18:20:54"
}

# de-de sheet: set Priority (column E) and Latest Handoff Datetime (column H)
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-17 18:21:01"
}
